$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.321.70"
$ws.Range("E2").Value = "  +3.34%  "
$ws.Range("D3").Value = "2.501.76"
$ws.Range("E3").Value = "  +2.57%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.60%  "
$ws.Range("E7").Value = "  +1.75%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.11%  "
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("D15").Value = "2.893.89"
$ws.Range("E15").Value = "  +2.72%  "
$ws.Range("D16").Value = "2.499.96"
$ws.Range("E16").Value = "  +2.80%  "
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("D18").Value = "47.277.17"
$ws.Range("E18").Value = "  +3.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.87%  "
$ws.Range("E20").Value = "  +4.92%  "
$ws.Range("D21").Value = "0.0₃0947"
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "250.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  +4.97%  "
$ws.Range("E29").Value = "  +3.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.77%  "
$ws.Range("E31").Value = "  +5.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("E34").Value = "  +4.00%  "
$ws.Range("E35").Value = "  +4.62%  "
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.24%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.67%  "
$ws.Range("E40").Value = "  +1.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "122.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("E44").Value = "  +2.57%  "
$ws.Range("D45").Value = "1.989.55"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.60%  "
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("E50").Value = "  +9.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.20%  "
